$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Project Timeline")
$ws1.Activate()
$ws1.Range("B9").Select()
Write-Output "ok"
